# LOQ4065.xlsx content update
# Rewrites the data grid on sheet1 to match the new course-syllabus content,
# re-sequences the "Avaliação" / "Bibliografia" rows, drops the old
# bibliography + detailed syllabus paragraphs, and removes the now-unused
# last row (old row 24), shrinking the used range from A1:C24 to A1:C23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Update cell values that change in place (style is already correct
#    because these cells already existed with the right formatting).
# ---------------------------------------------------------------------

$ws.Range("B10").Value = "5816812 - João Paulo Alves Silva"
$ws.Range("C10").Value = "5816812 - João Paulo Alves Silva"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("A14").Value = "Short syllabus:"
$shortSyllabus = @"
1 - Integrated Vision of Chemical Engineering. 2 - Multidisciplinary Project Studies aiming at integrating their knowledge in engineering. 3 - Development of multidisciplinary projects about Chemical Industry.
4 - Seminars: Presentation and discussion of the study results. 5. Final Report.
"@
$ws.Range("B14").Value = $shortSyllabus
$ws.Range("C14").Value = $shortSyllabus

$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2016"
$ws.Range("C15").Value = "01/01/2016"

$ws.Range("A18").Value = "Método:"

$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Provas escritas e Apresentação de Trabalhos"
$ws.Range("C19").Value = "Provas escritas e Apresentação de Trabalhos"

$ws.Range("A20").Value = "Norma de recuperação:"
$recNorm = "A nota será composta por ao menos uma prova escrita e trabalhos realizados e apresentados durante o semestre. O peso de cada atividade será definido segundo critérios do professor."
$ws.Range("B20").Value = $recNorm
$ws.Range("C20").Value = $recNorm

$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Média Final = (N + Prova Recuperação)/2"
$ws.Range("C21").Value = "Média Final = (N + Prova Recuperação)/2"

$ws.Range("A22").Value = "Requisitos:"

# ---------------------------------------------------------------------
# 2) Cells that are brand new in this layout (did not exist as <c> before)
#    need their number format / font / wrap copied from an untouched
#    donor cell in the same column so no extra cellXfs entries are
#    created and the style matches the rest of the column exactly.
# ---------------------------------------------------------------------

$ws.Range("A3").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A13").Value = "Programa resumido:"

$ws.Range("B3").Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null
$ws.Range("B18").Value = "5816812 - João Paulo Alves Silva"

$ws.Range("C3").Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$ws.Range("C18").Value = "5816812 - João Paulo Alves Silva"

$requisito = "LOQ4064 -  Engenharia de Processos Quimicos I  (Requisito fraco)`n"

$ws.Range("B3").Copy() | Out-Null
$ws.Range("B23").PasteSpecial(-4122) | Out-Null
$ws.Range("B23").Value = $requisito

$ws.Range("C3").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null
$ws.Range("C23").Value = $requisito

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Remove cells that no longer hold content in the new layout.
#    Using .Clear() (not just ClearContents) drops the <c> element
#    entirely instead of leaving an empty styled placeholder, and once
#    a whole row has no cells left, the <row> element itself disappears.
# ---------------------------------------------------------------------

$ws.Range("B17:C17").Clear()
$ws.Range("B22:C22").Clear()
$ws.Range("A23").Clear()
$ws.Range("B24:C24").Clear()

# ---------------------------------------------------------------------
# 4) Fix up row heights that changed between layouts.
# ---------------------------------------------------------------------

$ws.Rows(13).RowHeight = 60
$ws.Rows(15).RowHeight = 120
$ws.Rows(18).RowHeight = 60
$ws.Rows(21).RowHeight = 120
$ws.Rows(23).RowHeight = 30

# Rows 17 and 22 go back to the default (no explicit height) now that
# their long paragraphs are gone.
$ws.Rows(17).AutoFit()
$ws.Rows(22).AutoFit()

Write-Host "LOQ4065 content refreshed"
